$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40
$prev = 39

# Plain value cells (no special style in the source rows)
$ws.Cells.Item($row, 2).Value = "armenia"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 6).Value = "Noah"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Pyunik Yerevan"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 4.19
$ws.Cells.Item($row, 11).Value = "20/09/2023 04:12"
$ws.Cells.Item($row, 12).Value = 3.7
$ws.Cells.Item($row, 13).Value = "21/09/2023 16:56"
$ws.Cells.Item($row, 14).Value = 3.92
$ws.Cells.Item($row, 15).Value = "20/09/2023 04:12"
$ws.Cells.Item($row, 16).Value = 3.97
$ws.Cells.Item($row, 17).Value = "21/09/2023 16:56"
$ws.Cells.Item($row, 18).Value = 1.67
$ws.Cells.Item($row, 19).Value = "20/09/2023 04:12"
$ws.Cells.Item($row, 20).Value = 1.87
$ws.Cells.Item($row, 21).Value = "21/09/2023 16:56"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/noah-pyunik-yerevan/6FrxWTpA/"

# Column A carries the bold/centered/bordered "Indice" style used throughout the table
$ws.Cells.Item($row, 1).Value = 39
$ws.Cells.Item($prev, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

# Column E carries the date/time number format style
$ws.Cells.Item($row, 5).Value = 45190.70833333334
$ws.Cells.Item($prev, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$excel.CutCopyMode = $false
